$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(4, "Televisor LG", "Electrónica", 1000, 1500),
    @(5, "Televisor Samsung", "Electrónica", 1000, 1500),
    @(6, "Set vasos x 6 UND", "Hogar", 1000, 1500),
    @(7, "Silla Rimax", "Hogar", 1000, 1500),
    @(8, "Manzana", "Alimentos", 1000, 1500),
    @(9, "Granadilla", "Alimentos", 1000, 1500),
    @(10, "Camiseta Nike", "Ropa", 1000, 1500)
)

$row = 5
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $item[4]
    $row++
}
